$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet's display name (date rolled from 10-08 to 10-09)
$ws.Name = "Through 2021-10-09"

# Update the "October (through 10-08)" label to "October (through 10-09)"
$ws.Range("A11").Value = "October (through 10-09)"

# September row: H10 178 -> 179
$ws.Cells.Item(10, 8).Value = 179

# October row (row 11): B11:H11
$ws.Cells.Item(11, 2).Value = 7
$ws.Cells.Item(11, 3).Value = 16
$ws.Cells.Item(11, 4).Value = 17
$ws.Cells.Item(11, 5).Value = 23
$ws.Cells.Item(11, 6).Value = 8
$ws.Cells.Item(11, 7).Value = 38
$ws.Cells.Item(11, 8).Value = 60

# Total row (row 12): B12:H12
$ws.Cells.Item(12, 2).Value = 233
$ws.Cells.Item(12, 3).Value = 445
$ws.Cells.Item(12, 4).Value = 644
$ws.Cells.Item(12, 5).Value = 571
$ws.Cells.Item(12, 6).Value = 430
$ws.Cells.Item(12, 7).Value = 939
$ws.Cells.Item(12, 8).Value = 1310
